# Weekly price update for "Vega Modelo de Temuco - Chirimoya".
# Three new rows (for date 44461) are inserted before the previous last
# week's rows, pushing the existing rows 54-58 down to rows 57-61, and a
# few of those older rows get corrected values as part of the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three fresh rows above current row 54 (shifts 54-58 -> 57-61).
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()

# 2) Populate the three new rows (54-56) with this week's data.
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44461
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100107
$ws.Cells.Item(54, 8).Value = "Otros"
$ws.Cells.Item(54, 9).Value = 100107002
$ws.Cells.Item(54, 10).Value = "Chirimoya"
$ws.Cells.Item(54, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(54, 12).Value = "Especial"
$ws.Cells.Item(54, 13).Value = 40
$ws.Cells.Item(54, 14).Value = 3500
$ws.Cells.Item(54, 15).Value = 3500
$ws.Cells.Item(54, 16).Value = 3500
$ws.Cells.Item(54, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(54, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 19).Value = 3500
$ws.Cells.Item(54, 20).Value = 1

$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44461
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100107
$ws.Cells.Item(55, 8).Value = "Otros"
$ws.Cells.Item(55, 9).Value = 100107002
$ws.Cells.Item(55, 10).Value = "Chirimoya"
$ws.Cells.Item(55, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 40
$ws.Cells.Item(55, 14).Value = 30000
$ws.Cells.Item(55, 15).Value = 30000
$ws.Cells.Item(55, 16).Value = 30000
$ws.Cells.Item(55, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(55, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(55, 19).Value = 3750
$ws.Cells.Item(55, 20).Value = 8

$ws.Cells.Item(56, 1).Value = 10
$ws.Cells.Item(56, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value = "La Araucanía"
$ws.Cells.Item(56, 4).Value = 44461
$ws.Cells.Item(56, 5).Value = 9
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100107
$ws.Cells.Item(56, 8).Value = "Otros"
$ws.Cells.Item(56, 9).Value = 100107002
$ws.Cells.Item(56, 10).Value = "Chirimoya"
$ws.Cells.Item(56, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(56, 12).Value = "Segunda"
$ws.Cells.Item(56, 13).Value = 30
$ws.Cells.Item(56, 14).Value = 28000
$ws.Cells.Item(56, 15).Value = 28000
$ws.Cells.Item(56, 16).Value = 28000
$ws.Cells.Item(56, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(56, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 19).Value = 3500
$ws.Cells.Item(56, 20).Value = 8

# 3) Fix up values on the (now shifted) old rows 57 and 58.
$ws.Cells.Item(57, 4).Value = 44162
$ws.Cells.Item(57, 13).Value = 85
$ws.Cells.Item(57, 14).Value = 2200
$ws.Cells.Item(57, 15).Value = 2300
$ws.Cells.Item(57, 16).Value = 2247
$ws.Cells.Item(57, 19).Value = 2247

$ws.Cells.Item(58, 4).Value = 44411
$ws.Cells.Item(58, 13).Value = 10

# 4) Append three more rows (59-61) recreating the rows that used to be
#    54, 56 and 57 before the correction (duplicated at the bottom).
$ws.Cells.Item(59, 1).Value = 10
$ws.Cells.Item(59, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(59, 3).Value = "La Araucanía"
$ws.Cells.Item(59, 4).Value = 44425
$ws.Cells.Item(59, 5).Value = 9
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100107
$ws.Cells.Item(59, 8).Value = "Otros"
$ws.Cells.Item(59, 9).Value = 100107002
$ws.Cells.Item(59, 10).Value = "Chirimoya"
$ws.Cells.Item(59, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(59, 12).Value = "Especial"
$ws.Cells.Item(59, 13).Value = 35
$ws.Cells.Item(59, 14).Value = 4500
$ws.Cells.Item(59, 15).Value = 4500
$ws.Cells.Item(59, 16).Value = 4500
$ws.Cells.Item(59, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(59, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(59, 19).Value = 4500
$ws.Cells.Item(59, 20).Value = 1

$ws.Cells.Item(60, 1).Value = 10
$ws.Cells.Item(60, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(60, 3).Value = "La Araucanía"
$ws.Cells.Item(60, 4).Value = 44425
$ws.Cells.Item(60, 5).Value = 9
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100107
$ws.Cells.Item(60, 8).Value = "Otros"
$ws.Cells.Item(60, 9).Value = 100107002
$ws.Cells.Item(60, 10).Value = "Chirimoya"
$ws.Cells.Item(60, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 20
$ws.Cells.Item(60, 14).Value = 3500
$ws.Cells.Item(60, 15).Value = 3500
$ws.Cells.Item(60, 16).Value = 3500
$ws.Cells.Item(60, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(60, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(60, 19).Value = 3500
$ws.Cells.Item(60, 20).Value = 1

$ws.Cells.Item(61, 1).Value = 10
$ws.Cells.Item(61, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(61, 3).Value = "La Araucanía"
$ws.Cells.Item(61, 4).Value = 44425
$ws.Cells.Item(61, 5).Value = 9
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100107
$ws.Cells.Item(61, 8).Value = "Otros"
$ws.Cells.Item(61, 9).Value = 100107002
$ws.Cells.Item(61, 10).Value = "Chirimoya"
$ws.Cells.Item(61, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(61, 12).Value = "Segunda"
$ws.Cells.Item(61, 13).Value = 25
$ws.Cells.Item(61, 14).Value = 3000
$ws.Cells.Item(61, 15).Value = 3000
$ws.Cells.Item(61, 16).Value = 3000
$ws.Cells.Item(61, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(61, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(61, 19).Value = 3000
$ws.Cells.Item(61, 20).Value = 1

Write-Output "Applied Chirimoya weekly update"
